$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 00:37"

# --- Swap country labels whose global ranking changed (column A) ---
# Colombia overtook Sudafrica
$ws.Range("A9").Value = "Colombia"
$ws.Range("A10").Value = "Sudafrica"
# Japon overtook Polonia
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Polonia"
# Nigeria overtook Etiopia
$ws.Range("A53").Value = "Nigeria"
$ws.Range("A54").Value = "Etiopia"
# Cuba overtook Mozambique and Surinam
$ws.Range("A118").Value = "Cuba"
$ws.Range("A119").Value = "Mozambique"
$ws.Range("A120").Value = "Surinam"

# --- Refresh the daily statistics (columns B-H) ---
# Row 4
$ws.Range("B4").Value = 6295733
$ws.Range("C4").Value = 38162
$ws.Range("D4").Value = 3535914
$ws.Range("E4").Value = 2569927
$ws.Range("G4").Value = 992
$ws.Range("H4").Value = 189892

# Row 5
$ws.Range("D5").Value = 3210405
$ws.Range("E5").Value = 663680

# Row 9
$ws.Range("B9").Value = 633339
$ws.Range("C9").Value = 9270
$ws.Range("D9").Value = 479568
$ws.Range("E9").Value = 133423
$ws.Range("G9").Value = 296
$ws.Range("H9").Value = 20348

# Row 10
$ws.Range("B10").Value = 630595
$ws.Range("C10").Value = 2336
$ws.Range("D10").Value = 553456
$ws.Range("E10").Value = 62750
$ws.Range("G10").Value = 126
$ws.Range("H10").Value = 14389

# Row 23
$ws.Range("B23").Value = 247391
$ws.Range("C23").Value = 1390
$ws.Range("D23").Value = 223100
$ws.Range("E23").Value = 14898

# Row 27
$ws.Range("B27").Value = 129923
$ws.Range("C27").Value = 498
$ws.Range("D27").Value = 115050
$ws.Range("E27").Value = 5738
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 9135

# Row 34
$ws.Range("B34").Value = 99280
$ws.Range("C34").Value = 165
$ws.Range("D34").Value = 74626
$ws.Range("E34").Value = 19193
$ws.Range("G34").Value = 21
$ws.Range("H34").Value = 5461

# Row 43
$ws.Range("B43").Value = 75644
$ws.Range("C43").Value = 751
$ws.Range("D43").Value = 63688
$ws.Range("E43").Value = 9166
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 2790

# Row 47
$ws.Range("B47").Value = 69001
$ws.Range("C47").Value = 609
$ws.Range("D47").Value = 58428
$ws.Range("E47").Value = 9266
$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 1307

# Row 48
$ws.Range("B48").Value = 68517
$ws.Range("C48").Value = 595
$ws.Range("D48").Value = 47865
$ws.Range("E48").Value = 18574
$ws.Range("G48").Value = 20
$ws.Range("H48").Value = 2078

# Row 53
$ws.Range("B53").Value = 54463
$ws.Range("C53").Value = 216
$ws.Range("D53").Value = 42439
$ws.Range("E53").Value = 10997
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 1027

# Row 54
$ws.Range("B54").Value = 54409
$ws.Range("C54").Value = 1105
$ws.Range("D54").Value = 19903
$ws.Range("E54").Value = 33660
$ws.Range("G54").Value = 18
$ws.Range("H54").Value = 846

# Row 78
$ws.Range("B78").Value = 19460
$ws.Range("C78").Value = 51
$ws.Range("E78").Value = 1394
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 415

# Row 83
$ws.Range("B83").Value = 16617
$ws.Range("C83").Value = 163
$ws.Range("D83").Value = 11760
$ws.Range("E83").Value = 4209
$ws.Range("G83").Value = 6
$ws.Range("H83").Value = 648

# Row 90
$ws.Range("B90").Value = 11034
$ws.Range("C90").Value = 163
$ws.Range("E90").Value = 1422

# Row 118
$ws.Range("B118").Value = 4126
$ws.Range("C118").Value = 61
$ws.Range("D118").Value = 3458
$ws.Range("E118").Value = 570
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 98

# Row 119
$ws.Range("B119").Value = 4117
$ws.Range("C119").Value = 78
$ws.Range("D119").Value = 2170
$ws.Range("E119").Value = 1922
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 25

# Row 120
$ws.Range("B120").Value = 4089
$ws.Range("D120").Value = 3171
$ws.Range("E120").Value = 846
$ws.Range("H120").Value = 72

# Row 134
$ws.Range("B134").Value = 2777
$ws.Range("C134").Value = 48
$ws.Range("D134").Value = 1115
$ws.Range("E134").Value = 1550
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 112

# Row 145
$ws.Range("B145").Value = 1976
$ws.Range("C145").Value = 14
$ws.Range("D145").Value = 1175
$ws.Range("E145").Value = 230
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 571

# Row 147
$ws.Range("B147").Value = 1920
$ws.Range("C147").Value = 123
$ws.Range("E147").Value = 1201

# Row 154
$ws.Range("B154").Value = 1434
$ws.Range("C154").Value = 18
$ws.Range("D154").Value = 1055
$ws.Range("E154").Value = 349
$ws.Range("G154").Value = 2
$ws.Range("H154").Value = 30

# Row 156
$ws.Range("B156").Value = 1382
$ws.Range("C156").Value = 9
$ws.Range("D156").Value = 779
$ws.Range("E156").Value = 562

# Row 165
$ws.Range("D165").Value = 904
$ws.Range("E165").Value = 36
